$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above the old "Dabney things" row (row 26),
# pushing the trailing notes (rows 26-28) down to rows 29-31.
$ws.Rows("24:26").Insert()

# Fill in the new reference row (row 24) with the Knuth citation.
$ws.Range("B24").Value = "The Art of Computer Programing"
$ws.Range("C24").Value = "Donald Knuth"
$ws.Range("D24").Value = 1938
$ws.Range("E24").Value = "https://doc.lagout.org/science/0_Computer%20Science/2_Algorithms/The%20Art%20of%20Computer%20Programming%20%28vol.%203_%20Sorting%20and%20Searching%29%20%282nd%20ed.%29%20%5BKnuth%201998-05-04%5D.pdf"
$ws.Range("F24").Value = "Multibrackets"
$ws.Range("G24").Value = "x"
$ws.Range("H24").Value = "networked"

# Update what used to be row 28 ("Art of computing"), now shifted to row 31,
# replacing its single note with a link + a new note.
$ws.Range("C31").Value = "https://www-degruyter-com.ezp-prod1.hul.harvard.edu/document/doi/10.1515/jqas-2012-0055/pdf?stream=true"
$ws.Range("D31").Value = "(double elim designs?)"

# Move the selection cursor to match the saved workbook's recorded selection.
$ws.Range("D32").Select()
